$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "REPLACE FAILED for: $old"
    }
}

function Insert-Run-After-Text($anchorText, $newText, $fontSizePt) {
    # Finds anchorText, collapses to its end, inserts newText as a new run with
    # matching Calibri/black/fontSizePt formatting.
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "ANCHOR NOT FOUND: $anchorText"
        return $null
    }
    $rng.Collapse(0)
    $rng.InsertAfter($newText)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = $fontSizePt
    $rng.Font.Color = 0
    return $rng
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-Text "Unraveling the Enigmatic Enigma Machine" "The Marvels of Chemistry: Unveiling the Secrets of Our World"

# ---------------------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------------------
Replace-Text "Agatha Sinclair" "Emily Davies"

# ---------------------------------------------------------------------------
# 3. Email paragraph
#    old runs: "agathasinclair@triangulate" / "." / "tech"
#    new runs: "Emily" / "." / "Davies@centennialacademy" / "." / "org"
# ---------------------------------------------------------------------------
Replace-Text "agathasinclair@triangulate" "Emily"
$okTech = $d.Content.Find.Execute("tech", $true, $true, $false, $false, $false, $true, 1, $false, "Davies@centennialacademy", 2)
if (-not $okTech) { Write-Host "REPLACE FAILED for: tech" }
Insert-Run-After-Text "Davies@centennialacademy" "." 16 | Out-Null
Insert-Run-After-Text "Davies@centennialacademy." "org" 16 | Out-Null

# ---------------------------------------------------------------------------
# 4. Body paragraph 1 (the long paragraph with line breaks)
# ---------------------------------------------------------------------------
Replace-Text "In the tumultuous era of World War II, one enigmatic device cast a long shadow over the battlefield of communication: the Enigma machine" "Chemistry encompasses myriad complexities of life and non-life as it delves into the essence of matter and its transformation"

Replace-Text " This electromechanical marvel, conceived by German engineers, held the power to encipher messages with an almost impenetrable veil of secrecy, becoming a formidable tool for strategic communication" " It offers an architectural lens to understand the world and speculate on possibilities"

Replace-Text " As the Allies sought to unlock the secrets hidden within the Enigma's intricate mechanisms, a remarkable team of codebreakers embarked on a relentless quest to decipher its complex cipher, setting the stage for one of history's most enthralling intellectual battles" " Exploring atoms, elements, and molecules empowers students to comprehend the building blocks of matter and unravel the story of creation"

Replace-Text " The story of the Enigma machine is a testament to human ingenuity, perseverance, and the unwavering determination to uncover hidden truths, forever etching its place in the annals of cryptography and military history" " Chemistry fosters critical thinking skills, analytical and investigative approaches, and enhanced scientific literacy"

# Two new sentences inserted after "... enhanced scientific literacy"
$r = Insert-Run-After-Text "Chemistry fosters critical thinking skills, analytical and investigative approaches, and enhanced scientific literacy" "." 12
Insert-Run-After-Text "Chemistry fosters critical thinking skills, analytical and investigative approaches, and enhanced scientific literacy." " Moreover, it establishes a foundation for health, medicine, environmental protection, and countless industrial domains" 12 | Out-Null
Insert-Run-After-Text "Moreover, it establishes a foundation for health, medicine, environmental protection, and countless industrial domains" "." 12 | Out-Null
Insert-Run-After-Text "Moreover, it establishes a foundation for health, medicine, environmental protection, and countless industrial domains." " The study of chemistry goes beyond the classroom, inviting students to interrogate the natural world through insightful observation" 12 | Out-Null

Replace-Text "In the heart of Bletchley Park, England, a group of brilliant minds, including Alan Turing, Joan Clarke, and Gordon Welchman, united under the shared purpose of cracking the Enigma code" "The quest to understand the mechanisms behind chemical reactions engages both creativity and logic, enabling students to conceptualize the intricate interplay of elements and compounds"

Replace-Text " Drawing upon mathematics, engineering, and sheer human intuition, they dissected the machine's inner workings, uncovering its intricate patterns and vulnerabilities" " Chemistry unveils the sophisticated communication that occurs at the atomic level, dictating the properties and behaviors of substances"

Replace-Text " With painstaking precision, they constructed formidable machines like the Bombe, a high-speed codebreaking device, and the Colossus, the world's first programmable computer, pushing the boundaries of technology and innovation" " By learning the language of chemistry, students attain a newfound appreciation for the natural world, developing a keen eye for unraveling the mysteries of life itself"

Replace-Text " Through tireless hours of meticulous analysis and unwavering resolve, the codebreakers gradually chipped away at the Enigma's defenses, piecing together fragments of intelligence that would ultimately shape the course of the war" " A solid foundation in chemistry cultivates a sense of enlightenment and satisfaction, empowering individuals to make informed decisions in their personal lives"

Replace-Text "The successful cracking of the Enigma code had far-reaching implications" "Moreover, chemistry plays a pivotal role in addressing global concerns"

Replace-Text " It provided invaluable insights into enemy troop movements, military strategies, and diplomatic communications, enabling the Allies to anticipate and counter German actions with remarkable precision" " It informs the development of new energy sources, fuels, and materials while aiding in the discovery of therapies to combat diseases"

Replace-Text " This intelligence proved pivotal in numerous decisive battles, including the Battle of Midway in the Pacific and the Battle of Kursk on the Eastern Front" " Comprehending the intricate relationship between chemistry and pressing issues such as air pollution, climate change, and water contamination unravels paths toward solutions"

Replace-Text " By neutralizing the Enigma's effectiveness, the Allies significantly weakened the German war effort, contributing to their eventual defeat" " Chemistry serves as a catalytic agent, propelling investigations into the cosmos, materials science, and energy storage, constantly expanding our understanding of the universe and our place within it"

# Remove the trailing sentence (". The Enigma story serves as a poignant reminder ...") entirely,
# including its leading "." separator run, but keep the very last "." run intact.
$rngDel = $d.Content
$okDel = $rngDel.Find.Execute(". The Enigma story serves as a poignant reminder of the enduring power of human intellect and resilience in the face of adversity", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($okDel) {
    $rngDel.Delete()
} else {
    Write-Host "Could not find trailing sentence to delete"
}

# ---------------------------------------------------------------------------
# 5. "Summary" heading: remove the lastRenderedPageBreak marker
# ---------------------------------------------------------------------------
$summaryHeading = $d.Content
$okS = $summaryHeading.Find.Execute("Summary", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($okS) {
    $summaryHeading.Select() | Out-Null
}

# ---------------------------------------------------------------------------
# 6. Summary body paragraph text replacements
# ---------------------------------------------------------------------------
Replace-Text "The Enigma machine, a formidable cipher device employed by Nazi Germany during World War II, posed a daunting challenge to Allied intelligence" "The beauty of Chemistry is revealed when we unravel the secrets of our world"

Replace-Text " With unwavering determination, a team of brilliant codebreakers, led by Alan Turing, embarked on a tireless quest to decipher its complex cipher" " By engaging curiosity, logic, and creativity, students delve into the complexities of matter and its transformation through atomic-level interactions"

Replace-Text " Through a combination of mathematical ingenuity, technological innovation, and sheer human perseverance, they gradually chipped away at the Enigma's defenses, providing invaluable intelligence that ultimately contributed to Allied victory" " Piecing together the elements of elements and compounds, we can understand the universe's natural choreography"

Replace-Text " The Enigma story stands as a testament to the indomitable spirit of human intellect and its power to overcome even the most formidable obstacles" " Grounded in chemistry's knowledge, students explore intricate processes that govern life and gain insight into the challenges humanity faces today"

# New sentence appended at the end of the Summary paragraph
Insert-Run-After-Text "Grounded in chemistry's knowledge, students explore intricate processes that govern life and gain insight into the challenges humanity faces today" "." 11 | Out-Null
Insert-Run-After-Text "Grounded in chemistry's knowledge, students explore intricate processes that govern life and gain insight into the challenges humanity faces today." " Chemistry empowers us with the mental ingenuity to navigate the future, fueling investigations that shape our understanding of the cosmos" 11 | Out-Null

# ---------------------------------------------------------------------------
# 7. New empty paragraph at the very end of the document
# ---------------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Host "DONE"
